$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$bigrams = @(
    @('(''reddit'', ''kotakuinact'')', 13),
    @('(''kotakuinact'', ''comment'')', 13),
    @('(''al'', ''gore'')', 13),
    @('(''polit'', ''statement'')', 12),
    @('(''video'', ''game'')', 11),
    @('(''sea'', ''level'')', 11),
    @('(''ice'', ''cap'')', 11),
    @('(''man'', ''make'')', 10),
    @('(''year'', ''ago'')', 9),
    @('(''carbon'', ''emiss'')', 8),
    @('(''hockey'', ''stick'')', 8),
    @('(''last'', ''year'')', 8),
    @('(''level'', ''rise'')', 7),
    @('(''co2'', ''emiss'')', 7),
    @('(''nasa'', ''gov'')', 7),
    @('(''specif'', ''heat'')', 7),
    @('(''black'', ''peopl'')', 7),
    @('(''low'', ''iq'')', 7),
    @('(''gather'', ''storm'')', 6),
    @('(''melt'', ''ice'')', 6),
    @('(''lord'', ''believ'')', 6),
    @('(''power'', ''plant'')', 6),
    @('(''chang'', ''polit'')', 6),
    @('(''pari'', ''agreement'')', 6),
    @('(''co2'', ''level'')', 6),
    @('(''bell'', ''curv'')', 6),
    @('(''coal'', ''oil'')', 6),
    @('(''publish'', ''report'')', 5),
    @('(''carbon'', ''pollut'')', 5),
    @('(''global'', ''catastroph'')', 5),
    @('(''mind'', ''worm'')', 5),
    @('(''climat'', ''model'')', 5),
    @('(''natur'', ''disast'')', 5),
    @('(''carbon'', ''dioxid'')', 5),
    @('(''realli'', ''want'')', 5),
    @('(''chang'', ''thing'')', 5),
    @('(''make'', ''sen'')', 5),
    @('(''human'', ''hive'')', 5),
    @('(''degre'', ''celsius'')', 5),
    @('(''giss'', ''nasa'')', 5),
    @('(''climat'', ''scienc'')', 5),
    @('(''chang'', ''real'')', 5),
    @('(''high'', ''co2'')', 5),
    @('(''green'', ''tech'')', 5),
    @('(''late'', ''game'')', 5),
    @('(''global'', ''climat'')', 5),
    @('(''solar'', ''panel'')', 5),
    @('(''year'', ''year'')', 5),
    @('(''grow'', ''wheat'')', 5),
    @('(''get'', ''grip'')', 5),
    @('(''greenhous'', ''effect'')', 5),
    @('(''use'', ''coal'')', 5),
    @('(''think'', ''peopl'')', 5),
    @('(''import'', ''peopl'')', 5),
    @('(''unit'', ''nation'')', 4),
    @('(''nation'', ''publish'')', 4),
    @('(''report'', ''last'')', 4),
    @('(''last'', ''month'')', 4),
    @('(''month'', ''warn'')', 4),
    @('(''warn'', ''drastic'')', 4),
    @('(''drastic'', ''refor'')', 4),
    @('(''refor'', ''effort'')', 4),
    @('(''effort'', ''reduct'')', 4),
    @('(''reduct'', ''carbon'')', 4),
    @('(''pollut'', ''planet'')', 4),
    @('(''planet'', ''soon'')', 4),
    @('(''soon'', ''face'')', 4),
    @('(''face'', ''irrever'')', 4),
    @('(''irrever'', ''global'')', 4),
    @('(''web'', ''archiv'')', 4),
    @('(''archiv'', ''web'')', 4),
    @('(''time'', ''sinc'')', 4),
    @('(''polit'', ''issu'')', 4),
    @('(''thing'', ''happen'')', 4),
    @('(''effect'', ''climat'')', 4),
    @('(''call'', ''power'')', 4),
    @('(''ocean'', ''rise'')', 4)
)

for ($i = 0; $i -lt $bigrams.Length; $i++) {
    $row = 5 + $i
    $ws.Range("B$row").Value = $bigrams[$i][0]
    $ws.Range("C$row").Value = [int]$bigrams[$i][1]
}